# "ultimo cambios input form" - update the purchasing export sheet:
#  - rename a few header columns (cotizacion/cliente/vendedor -> documento/autorizacion/proveedor)
#  - drop the "Aprobado" column (shift "Estado" left)
#  - replace the sample data row with a new row, and drop the extra 3rd row
#  - retune a few column widths

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra third data row entirely.
$ws.Rows.Item(3).Delete()

# Remove the "Aprobado" column (old column L); "Estado" (old column M) shifts into L.
$ws.Columns.Item(12).Delete()

# Header row renames.
$ws.Range("B1").Value = "Numero Doc."
$ws.Range("C1").Value = "Autorizacion"
$ws.Range("D1").Value = "Proveedor"

# New data row values. The source export keeps dates/currency as plain text
# (not real numbers/dates), so force those cells to Text before assigning -
# otherwise Excel auto-parses them into a date serial / number.
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "123-232-323232323"
$ws.Range("C2").Value = 12345678901

$ws.Range("D2").Value = "qweqwe"

$textCells = @("E2", "F2", "G2", "H2", "I2", "J2")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("E2").Value = "2023-08-22"
$ws.Range("F2").Value = "$8610.00"
$ws.Range("G2").Value = "$0.00"
$ws.Range("H2").Value = "$1033.20"
$ws.Range("I2").Value = "$9643.20"
$ws.Range("J2").Value = "$9643.20"

$ws.Range("K2").Value = "Pagado"
$ws.Range("L2").Value = "Activo"

# Column width retune (ColumnWidth is in characters; the engine adds the
# standard ~0.8333 padding when round-tripping to the OOXML `width`, so we
# subtract it here to land on the exact target widths).
$offset = 5/6
$ws.Columns.Item(2).ColumnWidth = 17 - $offset
$ws.Columns.Item(3).ColumnWidth = 13 - $offset
$ws.Columns.Item(4).ColumnWidth = 35 - $offset
$ws.Columns.Item(5).ColumnWidth = 10 - $offset
$ws.Columns.Item(6).ColumnWidth = 14 - $offset
